$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 14, shifting existing rows 14-48 down to 15-49.
$ws.Rows.Item(14).Insert()

# Populate the new row 14 with the new weekly price data.
$ws.Cells.Item(14, 1).Value = 8
$ws.Cells.Item(14, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(14, 3).Value = "Coquimbo"
$ws.Cells.Item(14, 4).Value = 44497
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(14, 6).Value = 100112052
$ws.Cells.Item(14, 7).Value = "Albahaca"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 740
$ws.Cells.Item(14, 11).Value = 3800
$ws.Cells.Item(14, 12).Value = 4000
$ws.Cells.Item(14, 13).Value = 3900
$ws.Cells.Item(14, 14).Value = "$/paquete"
$ws.Cells.Item(14, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(14, 16).Value = 3900
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"
